# Applies the update to the Turkey Super Lig 2023-2024 results workbook:
#  - Rows 10/11, 18/19 and 46/47 had their match data (columns F:V) swapped
#    with each other (index/country/tournament/season/date in A:E unchanged).
#  - 8 new match rows (49-56) were appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$stage = 1000   # scratch row, far away from real data, cleared at the end

function Swap-RowData {
    param($r1, $r2)

    $ws.Range("F$r1`:V$r1").Copy()
    $ws.Range("F$stage`:V$stage").PasteSpecial(-4163)

    $ws.Range("F$r2`:V$r2").Copy()
    $ws.Range("F$r1`:V$r1").PasteSpecial(-4163)

    $ws.Range("F$stage`:V$stage").Copy()
    $ws.Range("F$r2`:V$r2").PasteSpecial(-4163)

    $ws.Range("F$stage`:V$stage").ClearContents()
}

Swap-RowData 10 11
Swap-RowData 18 19
Swap-RowData 46 47

# Append the 8 new rows (49-56), copying formatting from the last existing
# data row (48) so styles (bold/border index column, date format, etc.)
# match the rest of the sheet.
$ws.Range("A48:V48").Copy()
$ws.Range("A49:V56").PasteSpecial(-4122)

$newRows = @(
    @{ A=48; F="Pendikspor";   G=1; H="Karagumruk";      I=1; E=45192.66666666666;
       J=2.78; K="21/09/2023 09:10"; L=2.44; M="23/09/2023 15:59";
       N=3.49; O="21/09/2023 09:10"; P=3.33; Q="23/09/2023 15:59";
       R=2.55; S="21/09/2023 09:10"; T=3.15; U="23/09/2023 15:59";
       V="https://www.betexplorer.com/football/turkey/super-lig/pendikspor-f-karagumruk/U15E8uU6/" },

    @{ A=49; F="Rizespor";     G=1; H="Sivasspor";        I=1; E=45192.66666666666;
       J=2.16; K="19/09/2023 14:42"; L=1.96; M="23/09/2023 15:59";
       N=3.65; O="19/09/2023 14:42"; P=3.7;  Q="23/09/2023 15:57";
       R=3.39; S="19/09/2023 14:42"; T=4.01; U="23/09/2023 15:59";
       V="https://www.betexplorer.com/football/turkey/super-lig/rizespor-sivasspor/ALCo7mAa/" },

    @{ A=50; F="Antalyaspor";  G=2; H="Samsunspor";       I=0; E=45192.79166666666;
       J=1.97; K="19/09/2023 14:42"; L=2.18; M="23/09/2023 18:52";
       N=3.8;  O="19/09/2023 14:42"; P=3.55; Q="23/09/2023 18:52";
       R=3.82; S="19/09/2023 14:42"; T=3.47; U="23/09/2023 18:52";
       V="https://www.betexplorer.com/football/turkey/super-lig/antalyaspor-samsunspor/xWKY9kQt/" },

    @{ A=51; F="Basaksehir";   G=1; H="Galatasaray";      I=2; E=45192.79166666666;
       J=5.38; K="17/09/2023 15:12"; L=4.24; M="23/09/2023 18:57";
       N=4.5;  O="17/09/2023 15:12"; P=3.75; Q="23/09/2023 18:53";
       R=1.6;  S="17/09/2023 15:12"; T=1.9;  U="23/09/2023 18:57";
       V="https://www.betexplorer.com/football/turkey/super-lig/basaksehir-galatasaray/KjY6EVHP/" },

    @{ A=52; F="Ankaragucu";   G=1; H="Konyaspor";        I=1; E=45193.66666666666;
       J=2.23; K="18/09/2023 18:13"; L=2.51; M="24/09/2023 15:56";
       N=3.7;  O="18/09/2023 18:13"; P=3.43; Q="24/09/2023 15:41";
       R=3.13; S="18/09/2023 18:13"; T=2.97; U="24/09/2023 15:56";
       V="https://www.betexplorer.com/football/turkey/super-lig/ankaragucu-konyaspor/WUDs8Tfg/" },

    @{ A=53; F="Besiktas";     G=2; H="Kayserispor";      I=1; E=45193.66666666666;
       J=1.38; K="17/09/2023 18:12"; L=1.37; M="24/09/2023 15:30";
       N=5.56; O="17/09/2023 18:12"; P=5.4;  Q="24/09/2023 15:53";
       R=7.64; S="17/09/2023 18:12"; T=8.720000000000001; U="24/09/2023 15:53";
       V="https://www.betexplorer.com/football/turkey/super-lig/besiktas-kayserispor/l6Ox99um/" },

    @{ A=54; F="Alanyaspor";   G=0; H="Fenerbahce";       I=1; E=45193.79166666666;
       J=5.29; K="17/09/2023 18:12"; L=5.65; M="24/09/2023 18:59";
       N=4.73; O="17/09/2023 18:12"; P=4.49; Q="24/09/2023 18:59";
       R=1.57; S="17/09/2023 18:12"; T=1.58; U="24/09/2023 18:59";
       V="https://www.betexplorer.com/football/turkey/super-lig/alanyaspor-fenerbahce/2F3M6JbJ/" },

    @{ A=55; F="Kasimpasa";    G=2; H="Adana Demirspor";  I=1; E=45193.79166666666;
       J=2.79; K="17/09/2023 18:12"; L=3.51; M="24/09/2023 18:59";
       N=3.85; O="17/09/2023 18:12"; P=4.14; Q="24/09/2023 18:59";
       R=2.37; S="17/09/2023 18:12"; T=1.99; U="24/09/2023 18:58";
       V="https://www.betexplorer.com/football/turkey/super-lig/kasimpasa-adanademirspor/OO4I7aqD/" }
)

$r = 49
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = "turkey"
    $ws.Cells.Item($r, 3).Value = "super-lig"
    $ws.Cells.Item($r, 4).Value = "2023-2024"
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $ws.Cells.Item($r, 21).Value = $row.U
    $ws.Cells.Item($r, 22).Value = $row.V
    $r = $r + 1
}

Write-Host "done"
